$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Remove the "_GoBack" bookmark (w:bookmarkStart / w:bookmarkEnd)
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ------------------------------------------------------------------
# Helper: fill an empty table cell's (empty) paragraph with a run of
# text, applying the same character formatting that is already used
# by the surrounding cells (Segoe UI, 10.5pt / 21 half-points, black).
# ------------------------------------------------------------------
function Fill-Cell($table, $row, $col, $text) {
    $cell = $table.Cell($row, $col)
    $para = $cell.Range.Paragraphs.Item(1)
    $rng = $para.Range
    $rng.Collapse(1)
    $rng.InsertAfter($text)

    $cell2 = $table.Cell($row, $col)
    $para2 = $cell2.Range.Paragraphs.Item(1)
    $newRange = $d.Range($para2.Range.Start, $para2.Range.End - 1)
    $newRange.Font.Name = "Segoe UI"
    $newRange.Font.Size = 10.5
    $newRange.Font.Color = 0
}

# ------------------------------------------------------------------
# 2. Table 1 (Question 1 - Requirement / Reason table)
#    Row 2 is currently empty -> fill with the answer text.
# ------------------------------------------------------------------
$t1 = $d.Tables.Item(1)
Fill-Cell $t1 2 1 "Data Analysis and Predictive Data Analysis"
Fill-Cell $t1 2 2 "Databricks can be used to digest Big Data and perform analytics tasks"

# ------------------------------------------------------------------
# 3. Table 2 (Question 2 - Candidate data source / Reason table)
#    Row 2 is currently empty -> fill with the answer text.
# ------------------------------------------------------------------
$t2 = $d.Tables.Item(2)
Fill-Cell $t2 2 1 "Data Lake Store Gen II"
Fill-Cell $t2 2 2 "It can be used to store multiple extensions of files and data that can be accessed using Databricks"
